# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.438.71"
$ws.Range("E2").Value = "  +12.15%  "

$ws.Range("D3").Value = "1.826.31"
$ws.Range("E3").Value = "  +7.85%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.545"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.92%  "

$ws.Range("E7").Value = "  +0.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.40"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.73"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.25%  "

$ws.Range("E10").Value = "  +5.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0677"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.61%  "

$ws.Range("E12").Value = "  +3.40%  "

$ws.Range("D13").Value = "2.088.29"
$ws.Range("E13").Value = "  +7.92%  "

$ws.Range("D14").Value = "1.834.53"
$ws.Range("E14").Value = "  +8.42%  "

$ws.Range("E15").Value = "  +4.04%  "

$ws.Range("D16").Value = "34.423.03"
$ws.Range("E16").Value = "  +12.13%  "

$ws.Range("E17").Value = "  -5.26%  "

$ws.Range("E18").Value = "  +7.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "260.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.69%  "

$ws.Range("D21").Value = "0.0₃0750"
$ws.Range("E21").Value = "  +3.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.98%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.55%  "

$ws.Range("E25").Value = "  +0.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.15%  "

$ws.Range("E27").Value = "  +5.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.116"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("E31").Value = "  +9.56%  "

$ws.Range("E32").Value = "  +2.87%  "

$ws.Range("E33").Value = "  +6.74%  "

$ws.Range("E34").Value = "  +8.20%  "

$ws.Range("D35").Value = "1.584.55"
$ws.Range("E35").Value = "  +4.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.12%  "

$ws.Range("E37").Value = "  +3.55%  "

$ws.Range("E38").Value = "  +4.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.632"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "85.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.59%  "

$ws.Range("E43").Value = "  +7.11%  "

$ws.Range("E44").Value = "  +6.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0522"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.91%  "

$ws.Range("E46").Value = "  +4.40%  "

$ws.Range("D47").Value = "1.978.40"
$ws.Range("E47").Value = "  +8.06%  "

$ws.Range("E48").Value = "  +5.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.23%  "

$ws.Range("E51").Value = "  +7.10%  "
